$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the Locator Name / Locator Value table with the updated row set
$data = @(
    ,('Locator Name', 'Locator Value ')
    ,('url1', 'https://cls.dmo.lhubsg.com/billing')
    ,('url2', 'https://cls.dmo.lhubsg.com/billing/billing-configuration/charge-items')
    ,('url3', 'https://cls.dmo.lhubsg.com/billing/billing-configuration/charge-items/create')
    ,('billConfigLoc', '//a[@href=''/billing/billing-configuration'']')
    ,('chargeItemPage', '//a[@href=''/billing/billing-configuration/charge-items'']')
    ,('spinLoader', '//div[contains(@class, ''_custom-loading-modal-container'')]')
    ,('createBtn', '//button[contains(@class, ''ant-btn-primary'')]//span[contains(text(), ''Create'')]')
    ,('chargeItemName', 'charge-item-create_name')
    ,('chargeTypeId', 'charge-item-create_chargeTypeId')
    ,('registrationCTId', '//div[@title=''Registration'']')
    ,('adminCTId', '//div[@title=''Admin'']')
    ,('duplicateInvoiceId', '//div[@title=''Duplicate Invoice'']')
    ,('generalId', '//div[@title=''General'']')
    ,('chargeItemDesc', 'charge-item-create_description')
    ,('baseCurrencyId', 'charge-item-create_baseCurrencyId')
    ,('sgdBCId', '//div[@title=''SGD'']')
    ,('usdBCId', '//div[@title=''USD'']')
    ,('gstApplicableTrue', '//div[@id=''charge-item-create_gstApplicable'']//input[@class=''ant-radio-input'' and @value=''true'']')
    ,('gstApplicableFalse', '//div[@id=''charge-item-create_gstApplicable'']//input[@class=''ant-radio-input'' and @value=''false'']')
    ,('chargeAmount', 'charge-item-create_amount')
    ,('glNumberId', 'charge-item-create_glNumberId')
    ,('glNumberIdValue1', '//div[@title=''001'']')
    ,('glNumberIdValue2', '//div[@title=''002'']')
    ,('glNumberIdValue3', '//div[@title=''005'']')
    ,('glNumberIdValue4', '//div[@title=''006'']')
    ,('refundableId', 'charge-item-create_refundable')
    ,('tagToRegistrationTrue', '//div[@id=''charge-item-create_tagToRegistration'']//input[@class=''ant-radio-input'' and @value=''true'']')
    ,('tagToRegistrationFalse', '//div[@id=''charge-item-create_tagToRegistration'']//input[@class=''ant-radio-input'' and @value=''false'']')
    ,('saveBtn', '//button[@type=''button'']//span[contains(text(), ''Save'')]')
    ,('cancelBtn', '//button[contains(@class, ''ant-btn-default'')]//span[contains(text(), ''Cancel'')]')
    ,('xBtn', '//button[@aria-label=''Close'']')
    ,('funnelBtn', '//button[@data-testid=''filter-funnel-btn'']')
    ,('colExpBtn', '//div[@role=''button'']')
    ,('filterApplyBtn', '//button[contains(@class, ''ant-btn-primary'')]//span[contains(text(), ''Apply'')]')
    ,('modalOkBtn', '//button[contains(@class, ''confirm-modal-ok-btn'')]')
    ,('modalCancelBtn', '(//button[contains(@class, ''ant-btn-default'')]//span[contains(text(), ''Cancel'')])[2]')
    ,('modalConfirmContent', '//div[@class=''ant-modal-confirm-content'']')
    ,('searchCI', '//input[@placeholder=''Search by charge item'']')
    ,('searchCItem', '//td[contains(@class, ''ant-table-cell-ellipsis'')][1]')
    ,('profitCentreId', 'charge-item-create_profitCentreId')
    ,('profitCentreIdValue1', '//div[@title=''PC 1'']')
    ,('profitCentreIdValue2', '//div[@title=''PC 2'']')
    ,('errorMessage1', '//div[@id=''charge-item-create_name_help'']//div[@class=''ant-form-item-explain-error'']')
    ,('mandatoryErrorMessage1', '//div[@class=''ant-form-item-explain-error'']')
    ,('maxTextArea', '//span[@class=''ant-input-suffix'']//span[contains(text(), ''200/200'')]')
    ,('createDrawerTitle', '//h3[contains(text(), ''Create charge'')]')
    ,('chargeTotalItem', '//div[contains(text(), ''Total'')]')
    ,('statusFilterActive', '//input[@value=''ACTIVE'']')
    ,('statusFilterInactive', '//input[@value=''INACTIVE'']')
    ,('chargeItemRow', '//tr[contains(@class, ''ant-table-row'')]')
    ,('chargeItemColStatusList', '//td[contains(@class, ''ant-table-cell-ellipsis'')][9]')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Header row (row 1) becomes bold
$ws.Range("A1:B1").Font.Bold = $true

# Column A widens to fit the longest new locator name
$ws.Columns("A").ColumnWidth = 23.14

# Selection / active cell moves to B26 and the frozen top-left scroll resets
$ws.Range("B26").Select()
